$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-11-05 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-11-06 Monday", 2)

# Update table cell values (row index is 1-based in Word tables)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "12÷4="
$tbl.Cell(1,2).Range.Text = "59÷2="
$tbl.Cell(1,3).Range.Text = "58÷6="
$tbl.Cell(1,4).Range.Text = "51÷5="
$tbl.Cell(1,5).Range.Text = "81÷5="

$tbl.Cell(5,1).Range.Text = "52÷4="
$tbl.Cell(5,2).Range.Text = "22÷8="
$tbl.Cell(5,3).Range.Text = "34÷8="
$tbl.Cell(5,4).Range.Text = "67÷2="
$tbl.Cell(5,5).Range.Text = "24÷5="

$tbl.Cell(9,1).Range.Text = "16÷2="
$tbl.Cell(9,2).Range.Text = "55÷8="
$tbl.Cell(9,3).Range.Text = "74÷3="
$tbl.Cell(9,4).Range.Text = "26÷5="
$tbl.Cell(9,5).Range.Text = "32÷4="

$tbl.Cell(13,1).Range.Text = "59÷6="
$tbl.Cell(13,2).Range.Text = "30÷8="
$tbl.Cell(13,3).Range.Text = "99÷7="
$tbl.Cell(13,4).Range.Text = "89÷7="
$tbl.Cell(13,5).Range.Text = "18÷7="

$tbl.Cell(17,1).Range.Text = "22÷9="
$tbl.Cell(17,2).Range.Text = "56÷7="
$tbl.Cell(17,3).Range.Text = "15÷9="
$tbl.Cell(17,4).Range.Text = "50÷7="
$tbl.Cell(17,5).Range.Text = "27÷5="
